$d = $word.ActiveDocument

# The edit adds three new paragraphs ahead of the existing first paragraph
# ("Load windows 2008 server."): two lines of text followed by one blank
# line, all using the same "NoSpacing" style already used throughout the
# document. Collapse a range to the very start of the document and insert
# three fresh paragraph marks there (this naturally inherits the
# "NoSpacing" style from the following paragraph), then fill in the text
# for the first two of the newly created paragraphs.
$insertionPoint = $d.Paragraphs.Item(1).Range
$insertionPoint.Collapse(1)  # wdCollapseStart

$insertionPoint.InsertParagraphBefore() | Out-Null
$insertionPoint.InsertParagraphBefore() | Out-Null
$insertionPoint.InsertParagraphBefore() | Out-Null

$line1 = $d.Paragraphs.Item(1).Range
$line1.MoveEnd(1, -1) | Out-Null   # exclude the paragraph's own end mark
$line1.Text = "This is a test document for GIT"

$line2 = $d.Paragraphs.Item(2).Range
$line2.MoveEnd(1, -1) | Out-Null   # exclude the paragraph's own end mark
$line2.Text = "Need to check how changes are managed."

# Paragraph 3 is left empty, matching the blank line in the diff.
